$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H32").Value = 761.1667
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 775.8182
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 775.8182
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1427.8182

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

$ws.Range("H121").Value = 1018.5789
$ws.Range("I121").Value = 608.3333
$ws.Range("J121").Value = 1207.9231
$ws.Range("K121").Value = 1824.9999
$ws.Range("L121").Value = 3623.7693
$ws.Range("M121").Value = -77.99990000000003
$ws.Range("N121").Value = -7117.7693

$ws.Range("H138").Value = 2750.976
$ws.Range("I138").Value = 1284.3158
$ws.Range("J138").Value = 3962.5652
$ws.Range("K138").Value = 3852.9474
$ws.Range("L138").Value = 11887.6956
$ws.Range("M138").Value = 1287.0526
$ws.Range("N138").Value = -22167.6956

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H103").Value = 27533.334
$ws.Range("J103").Value = 27533.334
$ws.Range("L103").Value = 27533.334
$ws.Range("N103").Value = -29877.334

$ws.Range("H122").Value = 1038.6154
$ws.Range("I122").Value = 857.4286
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 2572.2858
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -122.2857999999997
$ws.Range("N122").Value = -8650

$ws.Range("H132").Value = 8214.9
$ws.Range("I132").Value = 6704.95
$ws.Range("J132").Value = 11234.8
$ws.Range("K132").Value = 20114.85
$ws.Range("L132").Value = 33704.39999999999
$ws.Range("M132").Value = -17584.85
$ws.Range("N132").Value = -38764.39999999999

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H132").Value = 40780
$ws.Range("J132").Value = 40780
$ws.Range("L132").Value = 40780
$ws.Range("N132").Value = -50900

$ws.Range("H134").Value = 761.22
$ws.Range("I134").Value = 678.8222
$ws.Range("J134").Value = 1502.8
$ws.Range("K134").Value = 2036.4666
$ws.Range("L134").Value = 4508.4
$ws.Range("M134").Value = 498.5334000000003
$ws.Range("N134").Value = -9578.4

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H58").Value = 1161.9512
$ws.Range("I58").Value = 804.0833
$ws.Range("J58").Value = 1667.1765
$ws.Range("K58").Value = 804.0833
$ws.Range("L58").Value = 1667.1765
$ws.Range("M58").Value = -601.0833
$ws.Range("N58").Value = -2073.1765

$ws.Range("H107").Value = 39253.42
$ws.Range("I107").Value = 63141.25
$ws.Range("J107").Value = 1032.9
$ws.Range("K107").Value = 63141.25
$ws.Range("L107").Value = 1032.9
$ws.Range("M107").Value = -61221.25
$ws.Range("N107").Value = -4872.9

$ws.Range("H134").Value = 1590.3903
$ws.Range("I134").Value = 1605.6552
$ws.Range("J134").Value = 1553.5
$ws.Range("K134").Value = 4816.9656
$ws.Range("L134").Value = 4660.5
$ws.Range("M134").Value = -2281.9656
$ws.Range("N134").Value = -9730.5

$ws.Range("H136").Value = 1161.9512
$ws.Range("I136").Value = 804.0833
$ws.Range("J136").Value = 1667.1765
$ws.Range("K136").Value = 2412.2499
$ws.Range("L136").Value = 5001.529500000001
$ws.Range("M136").Value = 137.7501000000002
$ws.Range("N136").Value = -10101.5295

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H8").Value = 338.7
$ws.Range("I8").Value = 338.7
$ws.Range("K8").Value = 1016.1
$ws.Range("M8").Value = -877.0999999999999

$ws.Range("H12").Value = 39.6875
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = 48.26923
$ws.Range("K12").Value = 7.5
$ws.Range("L12").Value = 144.80769
$ws.Range("M12").Value = 165.5
$ws.Range("N12").Value = -490.80769

$ws.Range("H23").Value = 83.28570999999999
$ws.Range("I23").Value = 86.2
$ws.Range("J23").Value = 76
$ws.Range("K23").Value = 258.6
$ws.Range("L23").Value = 228
$ws.Range("M23").Value = -23.60000000000002
$ws.Range("N23").Value = -698

$ws.Range("H49").Value = 1375
$ws.Range("I49").Value = 833.3333
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 2499.9999
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = -2343.9999
$ws.Range("N49").Value = -9312

$ws.Range("H57").Value = 8777.666999999999
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 8777.666999999999
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 26333.001
$ws.Range("M57").Value = ""
$ws.Range("N57").Value = -27451.001

$ws.Range("H113").Value = 363957.75
$ws.Range("I113").Value = 420.94287
$ws.Range("J113").Value = 788084
$ws.Range("K113").Value = 1262.82861
$ws.Range("L113").Value = 2364252
$ws.Range("M113").Value = 907.17139
$ws.Range("N113").Value = -2368592

$ws.Range("H132").Value = 948.7857
$ws.Range("I132").Value = 935.6667
$ws.Range("J132").Value = 958.625
$ws.Range("K132").Value = 8421.0003
$ws.Range("L132").Value = 8627.625
$ws.Range("M132").Value = -5891.0003
$ws.Range("N132").Value = -13687.625

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 28.566668
$ws.Range("I2").Value = 19.88889
$ws.Range("J2").Value = 41.583332
$ws.Range("K2").Value = 19.88889
$ws.Range("L2").Value = 41.583332
$ws.Range("M2").Value = 93.11111
$ws.Range("N2").Value = -267.583332

$ws.Range("H109").Value = 17713.75
$ws.Range("J109").Value = 17713.75
$ws.Range("L109").Value = 17713.75
$ws.Range("N109").Value = -19793.75

$ws.Range("H132").Value = 4684.718
$ws.Range("I132").Value = 5997.846
$ws.Range("J132").Value = 2058.4614
$ws.Range("K132").Value = 17993.538
$ws.Range("L132").Value = 6175.3842
$ws.Range("M132").Value = -15463.538
$ws.Range("N132").Value = -11235.3842

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H40").Value = 1856.0588
$ws.Range("I40").Value = 1576
$ws.Range("J40").Value = 2766.25
$ws.Range("K40").Value = 1576
$ws.Range("L40").Value = 2766.25
$ws.Range("M40").Value = -1440
$ws.Range("N40").Value = -3038.25

$ws.Range("H122").Value = 34733.547
$ws.Range("I122").Value = 43822.543
$ws.Range("K122").Value = 131467.629
$ws.Range("M122").Value = -129017.629

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = ""

$ws.Range("H130").Value = 34367.5
$ws.Range("J130").Value = 34367.5
$ws.Range("L130").Value = 34367.5
$ws.Range("N130").Value = -44407.5

$ws.Range("H131").Value = 46678
$ws.Range("J131").Value = 46678
$ws.Range("L131").Value = 46678
$ws.Range("N131").Value = -56758

$ws.Range("H132").Value = 2784.5593
$ws.Range("I132").Value = 4193.875
$ws.Range("J132").Value = 1114.2593
$ws.Range("K132").Value = 12581.625
$ws.Range("L132").Value = 3342.7779
$ws.Range("M132").Value = -10051.625
$ws.Range("N132").Value = -8402.777900000001

$ws.Range("H136").Value = 3027.4182
$ws.Range("I136").Value = 3374.9429
$ws.Range("J136").Value = 2419.25
$ws.Range("K136").Value = 10124.8287
$ws.Range("L136").Value = 7257.75
$ws.Range("M136").Value = -7574.8287
$ws.Range("N136").Value = -12357.75
